# Append a new trade record (row 9) to the BIIB random-trade sheet, copying
# the cell formatting used by the previous rows (date format in column A,
# boolean format in column G) so the new row matches the existing style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring the formatting of the last existing row down into the new row first,
# so the newly entered values inherit the same number formats / styles as
# the rows above them (date serial in A, IsShortSell boolean style in G).
$ws.Range("A8:I8").Copy()
$ws.Range("A9:I9").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(9, 1).Value = 42654.743703703702     # A9 Date
$ws.Cells.Item(9, 2).Value = $true                  # B9 Profitable
$ws.Cells.Item(9, 3).Value = 10402.91               # C9 Principle
$ws.Cells.Item(9, 4).Value = 10345.49               # D9 Start Principle
$ws.Cells.Item(9, 5).Value = 308                    # E9 BuyPrice
$ws.Cells.Item(9, 6).Value = 304.58999599999999     # F9 SellPrice
$ws.Cells.Item(9, 7).Value = $true                  # G9 IsShortSell
$ws.Cells.Item(9, 8).Value = -1.1100000000000001    # H9 Price Change %
$ws.Cells.Item(9, 9).Value = $false                 # I9 Strong trade

# The repeater that appended this row also nudged column A a little wider
# to keep the longer dates fitting ("repeater" width fix mentioned in the
# commit message).
$ws.Columns.Item(1).ColumnWidth = 14.541666666666666
